$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the "Price" column as text (source data has values like "25.784.16"
# that Excel would otherwise auto-convert to numbers), matching the original
# inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.784.16"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.627.65"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "215.47"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "0.5067"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "0.06418"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").Value = "0.07787"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "4.256"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "1.627.95"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "1.852.96"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "0.5576"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "63.09"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "0.0₅7539"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "25.809.32"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "193.83"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "4.313"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").Value = "9.808"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "5.996"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "1.787"
$ws.Range("E25").Value = "  -5.42%  "
$ws.Range("D26").Value = "141.15"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "0.1265"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "6.732"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "15.39"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").Value = "0.04874"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "3.276"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "3.186"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "1.554"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "2.374"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.8940"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "2.559"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "1.130.73"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "0.5467"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "0.9978"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "5.574"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "0.7957"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "97.42"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "1.781.90"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  -8.97%  "
$ws.Range("D47").Value = "0.4434"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").Value = "55.23"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").Value = "7.620"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "0.9999"
$ws.Range("E51").Value = "  -0.29%  "
